$d = $word.ActiveDocument

# Remove the 4 "SourceCode" style paragraphs (DEFINEDNAME verbatim blocks)
# that follow the date paragraph ("12 December, 2016") and precede the
# "Introduction" heading.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style -ne $null -and $p.Style.NameLocal -eq "Source Code") {
        $p.Range.Delete()
    }
}
